$d = $word.ActiveDocument

# Update the date/title paragraph
$d.Content.Find.Execute("2023-09-02 Saturday", $false, $false, $false, $false, $false, $true, 1, $false, "2023-09-03 Sunday", 2)

# Update the answer table by directly addressing cells (row, column),
# since several answer strings repeat across the table and a global
# Find/Replace would be ambiguous.
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text = "82÷5=16, 2"
$t.Cell(1, 2).Range.Text = "98÷8=12, 2"
$t.Cell(1, 3).Range.Text = "23÷2=11, 1"
$t.Cell(1, 4).Range.Text = "48÷8=6, 0"
$t.Cell(1, 5).Range.Text = "76÷7=10, 6"

$t.Cell(5, 1).Range.Text = "38÷4=9, 2"
$t.Cell(5, 2).Range.Text = "51÷2=25, 1"
$t.Cell(5, 3).Range.Text = "51÷8=6, 3"
$t.Cell(5, 4).Range.Text = "83÷4=20, 3"
$t.Cell(5, 5).Range.Text = "85÷9=9, 4"

$t.Cell(9, 1).Range.Text = "50÷4=12, 2"
$t.Cell(9, 2).Range.Text = "79÷9=8, 7"
$t.Cell(9, 3).Range.Text = "78÷8=9, 6"
$t.Cell(9, 4).Range.Text = "78÷4=19, 2"
$t.Cell(9, 5).Range.Text = "78÷8=9, 6"

$t.Cell(13, 1).Range.Text = "63÷5=12, 3"
$t.Cell(13, 2).Range.Text = "39÷5=7, 4"
$t.Cell(13, 3).Range.Text = "98÷2=49, 0"
$t.Cell(13, 4).Range.Text = "98÷7=14, 0"
$t.Cell(13, 5).Range.Text = "23÷5=4, 3"

$t.Cell(17, 1).Range.Text = "69÷2=34, 1"
$t.Cell(17, 2).Range.Text = "71÷5=14, 1"
$t.Cell(17, 3).Range.Text = "92÷8=11, 4"
$t.Cell(17, 4).Range.Text = "81÷7=11, 4"
$t.Cell(17, 5).Range.Text = "62÷9=6, 8"
